$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in H1, matching the style used by the other header cells (G1)
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill H2:H21 with 0 for Control rows (2-6, 12-16), 1 for MDD rows (7-11, 17-21)
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
for ($r = 7; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
for ($r = 12; $r -le 16; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
for ($r = 17; $r -le 21; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}

# Update refit numeric values (small precision changes from refitting NCDEs)
$ws.Range("D4").Value = 0.4924838018176934
$ws.Range("E4").Value = 0.4924838018176934

$ws.Range("D6").Value = [double]"7.401066327486236E-111"
$ws.Range("E6").Value = [double]"7.401066327486236E-111"

$ws.Range("D10").Value = 0.5129892445513994
$ws.Range("E10").Value = 0.4870107554486006

$ws.Range("D11").Value = [double]"1.372125161332212E-38"
$ws.Range("F11").Value = 9.275491714477539
